# Liability returns report: add an aggregated "Total" sheet that sums the
# Market Value across the IBT, Pension and Retirement plans, extend the
# Retirement sheet with a blank formatted row, and tidy up sheet selections.

$wb = $excel.ActiveWorkbook

$retirement = $wb.Worksheets.Item("Retirement")

# --- Add the new "Total" worksheet by duplicating Retirement (same column
# widths / styles / dimension as the source table) and placing it right
# after the Retirement tab. ---
$retirement.Copy($null, $retirement) | Out-Null
$total = $wb.Worksheets.Item($wb.Worksheets.Count)
$total.Name = "Total"

# Replace the copied Market Value numbers with a formula that aggregates
# the three plans for every dated row.
$total.Range("B2:B146").FormulaR1C1 = "=SUM(IBT!RC,Pension!RC,Retirement!RC)"

# Add the trailing formatted-but-empty row that mirrors the Retirement sheet.
$total.Range("B150").NumberFormat = """$""#,##0.00"

# Put the selection/active cell where the author left it on the new sheet.
$total.Activate() | Out-Null
$total.Range("D4").Select() | Out-Null

# --- Extend Retirement with the same blank formatted row and move its
# selection down from E3 to E7. ---
$retirement.Range("B150").NumberFormat = """$""#,##0.00"

$retirement.Activate() | Out-Null
$retirement.Range("E7").Select() | Out-Null
